$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1705882352941177
$ws.Range("C2").Value = 0.5911764705882353
$ws.Range("J2").Value = 0.03823529411764706
$ws.Range("P2").Value = 0.1264705882352941
$ws.Range("S2").Value = 0.07352941176470588
$ws.Range("B3").Value = 0.00966183574879227
$ws.Range("C3").Value = 0.02415458937198068
$ws.Range("J3").Value = 0.05314009661835749
$ws.Range("P3").Value = 0.6956521739130435
$ws.Range("S3").Value = 0.2173913043478261
$ws.Range("J4").Value = 0.07142857142857142
$ws.Range("P4").Value = 0.7142857142857143
$ws.Range("S4").Value = 0.2142857142857143
$ws.Range("J5").Value = 0.5
$ws.Range("P5").Value = 0.5
$ws.Range("B6").Value = 0.05531914893617021
$ws.Range("D6").Value = 0.01702127659574468
$ws.Range("F6").Value = 0.05957446808510639
$ws.Range("J6").Value = 0.2297872340425532
$ws.Range("O6").Value = 0.02553191489361702
$ws.Range("Q6").Value = 0.1574468085106383
$ws.Range("R6").Value = 0.09361702127659574
$ws.Range("S6").Value = 0.3617021276595745
$ws.Range("B7").Value = 0.09663865546218488
$ws.Range("D7").Value = 0.02941176470588235
$ws.Range("F7").Value = 0.05882352941176471
$ws.Range("J7").Value = 0.1848739495798319
$ws.Range("O7").Value = 0.02100840336134454
$ws.Range("Q7").Value = 0.1848739495798319
$ws.Range("R7").Value = 0.08823529411764706
$ws.Range("S7").Value = 0.3361344537815126
$ws.Range("B8").Value = 0.08097165991902834
$ws.Range("D8").Value = 0.02024291497975709
$ws.Range("E8").Value = 0.002024291497975709
$ws.Range("F8").Value = 0.04048582995951417
$ws.Range("J8").Value = 0.1518218623481781
$ws.Range("O8").Value = 0.02024291497975709
$ws.Range("Q8").Value = 0.1842105263157895
$ws.Range("R8").Value = 0.1052631578947368
$ws.Range("S8").Value = 0.3947368421052632
$ws.Range("B9").Value = 0.09956709956709957
$ws.Range("D9").Value = 0.01298701298701299
$ws.Range("F9").Value = 0.06926406926406926
$ws.Range("J9").Value = 0.1125541125541126
$ws.Range("O9").Value = 0.01731601731601732
$ws.Range("Q9").Value = 0.1688311688311688
$ws.Range("R9").Value = 0.1125541125541126
$ws.Range("S9").Value = 0.4069264069264069
$ws.Range("B10").Value = 0.1239495798319328
$ws.Range("D10").Value = 0.02310924369747899
$ws.Range("E10").Value = 0.0007002801120448179
$ws.Range("F10").Value = 0.07352941176470588
$ws.Range("J10").Value = 0.1400560224089636
$ws.Range("O10").Value = 0.01400560224089636
$ws.Range("Q10").Value = 0.2002801120448179
$ws.Range("R10").Value = 0.0819327731092437
$ws.Range("S10").Value = 0.342436974789916
$ws.Range("G11").Value = 0.1198830409356725
$ws.Range("J11").Value = 0.07602339181286549
$ws.Range("K11").Value = 0.1783625730994152
$ws.Range("L11").Value = 0.6052631578947368
$ws.Range("S11").Value = 0.02046783625730994
$ws.Range("G12").Value = 0.7725118483412322
$ws.Range("J12").Value = 0.1374407582938389
$ws.Range("K12").Value = 0.01895734597156398
$ws.Range("L12").Value = 0.02369668246445497
$ws.Range("S12").Value = 0.04739336492890995
$ws.Range("G13").Value = 0.7959183673469388
$ws.Range("J13").Value = 0.2040816326530612
$ws.Range("F15").Value = 0.03813559322033899
$ws.Range("H15").Value = 0.1694915254237288
$ws.Range("I15").Value = 0.09745762711864407
$ws.Range("J15").Value = 0.3050847457627119
$ws.Range("K15").Value = 0.07203389830508475
$ws.Range("M15").Value = 0.008474576271186441
$ws.Range("O15").Value = 0.05084745762711865
$ws.Range("S15").Value = 0.2584745762711864
$ws.Range("F16").Value = 0.0179372197309417
$ws.Range("H16").Value = 0.1883408071748879
$ws.Range("I16").Value = 0.05829596412556054
$ws.Range("J16").Value = 0.42152466367713
$ws.Range("K16").Value = 0.1121076233183857
$ws.Range("M16").Value = 0.0179372197309417
$ws.Range("O16").Value = 0.07174887892376682
$ws.Range("S16").Value = 0.1121076233183857
$ws.Range("F17").Value = 0.01221995926680244
$ws.Range("H17").Value = 0.175152749490835
$ws.Range("I17").Value = 0.1038696537678208
$ws.Range("J17").Value = 0.3727087576374745
$ws.Range("K17").Value = 0.1303462321792261
$ws.Range("M17").Value = 0.02036659877800407
$ws.Range("O17").Value = 0.05295315682281059
$ws.Range("S17").Value = 0.1323828920570265
$ws.Range("F18").Value = 0.0125
$ws.Range("H18").Value = 0.1958333333333333
$ws.Range("I18").Value = 0.08749999999999999
$ws.Range("J18").Value = 0.4083333333333333
$ws.Range("K18").Value = 0.1
$ws.Range("M18").Value = 0.01666666666666667
$ws.Range("O18").Value = 0.075
$ws.Range("S18").Value = 0.1041666666666667
$ws.Range("F19").Value = 0.01076040172166428
$ws.Range("H19").Value = 0.2058823529411765
$ws.Range("I19").Value = 0.08823529411764706
$ws.Range("J19").Value = 0.3651362984218077
$ws.Range("K19").Value = 0.09899569583931134
$ws.Range("M19").Value = 0.02152080344332855
$ws.Range("N19").Value = 0.001434720229555237
$ws.Range("O19").Value = 0.06527977044476327
$ws.Range("S19").Value = 0.1427546628407461
